# "added relative path to signup"
# Adds a new signup record (email / password / id) as row 11 of the
# credentials sheet, and switches the workbook's default font from
# Arial to Calibri (the font actually referenced by the theme's minor
# font scheme).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default/normal style font: Arial -> Calibri
$wb.Styles("Normal").Font.Name = "Calibri"

# New signup row
$ws.Range("A11").Value = "ritay12"
$ws.Range("B11").Value = "ritay12@"
$ws.Range("C11").Value = 311434621
